$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.624.04"
$ws.Range("E2").Value = "  -1.15%  "
$ws.Range("D3").Value = "3.445.95"
$ws.Range("E3").Value = "  -3.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.75"
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.75"
$ws.Range("E6").Value = "  -7.56%  "
$ws.Range("D7").Value = "3.445.52"
$ws.Range("E7").Value = "  -3.99%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("E11").Value = "  -9.98%  "
$ws.Range("E12").Value = "  -7.99%  "
$ws.Range("D13").Value = "4.035.08"
$ws.Range("E13").Value = "  -3.88%  "
$ws.Range("E14").Value = "  -11.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.54"
$ws.Range("E15").Value = "  -10.21%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "65.579.59"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.409.78"
$ws.Range("E17").Value = "  -4.98%  "
$ws.Range("E18").Value = "  -2.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.89"
$ws.Range("E19").Value = "  -10.91%  "
$ws.Range("E20").Value = "  -8.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.74"
$ws.Range("E21").Value = "  -7.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.14"
$ws.Range("E22").Value = "  -6.70%  "
$ws.Range("E23").Value = "  -10.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.25"
$ws.Range("E24").Value = "  -6.28%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").Value = "3.588.81"
$ws.Range("E26").Value = "  -3.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000105"
$ws.Range("E27").Value = "  -11.92%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("E29").Value = "  -11.10%  "
$ws.Range("E30").Value = "  -9.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.16"
$ws.Range("E31").Value = "  -12.67%  "
$ws.Range("D32").Value = "3.453.08"
$ws.Range("E32").Value = "  -3.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.145"
$ws.Range("E34").Value = "  -7.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.82"
$ws.Range("E35").Value = "  -8.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "173.19"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -13.59%  "
$ws.Range("E38").Value = "  -11.02%  "
$ws.Range("E39").Value = "  -8.74%  "
$ws.Range("E40").Value = "  -13.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0776"
$ws.Range("E41").Value = "  -8.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.815"
$ws.Range("E42").Value = "  -7.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.53"
$ws.Range("E43").Value = "  -5.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.45"
$ws.Range("E45").Value = "  -14.22%  "
$ws.Range("E46").Value = "  -12.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.02"
$ws.Range("E47").Value = "  -3.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.52"
$ws.Range("E49").Value = "  -8.50%  "
$ws.Range("E50").Value = "  -15.93%  "
$ws.Range("D51").Value = "2.210.11"
$ws.Range("E51").Value = "  -7.57%  "
